# Update Teradyne (TER) historical price rows: fix fixed_ticker (col I) to reference
# the company's own ticker "TER" and correct the open/close/high/low price and
# shares_outstanding values (cols D-H) that had been populated from other tickers.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each entry: row, open_price, close_price, high_price, low_price, shares_outstanding
$rowData = @(
    @(2, 17.52421361425563, 17.02058982849121, 18.11177291878725, 16.8154107558799, 159073532),
    @(3, 18.26608436530506, 18.01355743408203, 18.70566689513257, 16.90992280184338, 159073532),
    @(4, 16.89282459760514, 18.3193302154541, 19.10766217321126, 16.35788543841881, 159073532),
    @(5, 19.25753692694876, 18.28807258605957, 19.84109946803106, 17.44096717727977, 159073532),
    @(6, 20.22618603641685, 17.85607719421387, 20.62277799011307, 17.56335511919017, 159073532),
    @(7, 18.51736378414104, 18.70680046081543, 20.26017470577808, 17.8732830579726, 159073532),
    @(8, 20.51711592064236, 22.12239074707031, 22.49283808712042, 19.20629918985098, 159073532),
    @(9, 24.28901611807751, 27.02164840698243, 27.46915397823873, 24.03193868132056, 159073532),
    @(10, 29.74879420451959, 33.6618537902832, 34.92166768746762, 29.08071019461514, 159073532),
    @(11, 28.85114324144702, 33.07792663574219, 35.24869588332169, 28.3825634088631, 159073532),
    @(12, 35.90292294674281, 41.09624481201172, 41.81487815422734, 35.740030821825, 159073532),
    @(13, 41.95027299479072, 43.99452209472656, 46.57622374445597, 41.34563887627296, 159073532),
    @(14, 43.62297967224122, 31.30352592468262, 43.93072526295126, 30.99578033397258, 159073532),
    @(15, 36.11943933133284, 41.69110488891602, 44.14919241860265, 36.06160063659631, 159073532),
    @(16, 36.17920659571675, 33.281005859375, 36.3820797628, 29.19454142073706, 159073532),
    @(17, 29.57638445840699, 34.85442352294922, 36.28772523045483, 28.646673994642, 159073532),
    @(18, 38.94867833554941, 47.55756759643555, 48.07196459179146, 38.94867833554941, 159073532),
    @(19, 48.0490330281303, 54.2059211730957, 57.00716129796547, 42.59245213160131, 159073532),
    @(20, 57.23085330481264, 59.64713668823242, 63.75871761957458, 55.42838128268901, 159073532),
    @(21, 67.16640776841714, 64.38569641113281, 79.58692842283936, 63.98566790762168, 159073532),
    @(22, 51.19164483671422, 61.10947799682617, 65.99510936093237, 49.9800103457828, 159073532),
    @(23, 82.89468491864777, 87.06388854980469, 88.76679976436392, 80.97646163172745, 159073532),
    @(24, 79.20974844034899, 86.07838439941406, 91.29109910276996, 76.73077050936453, 159073532),
    @(25, 119.5324919979134, 111.2942886352539, 140.6379964654262, 110.8627609892769, 159073532),
    @(26, 121.7971843855428, 122.7590026855469, 136.6955028293391, 120.4133524041486, 159073532),
    @(27, 131.0763456091435, 124.7410736083984, 131.4888733454254, 114.3197911901227, 159073532),
    @(28, 107.836257594677, 135.8913879394531, 136.38289322931, 103.4520313091272, 159073532),
    @(29, 162.167749496305, 115.5123062133789, 164.9515358987887, 100.8359595156407, 159073532),
    @(30, 117.7260779982724, 103.8336715698242, 118.346359891589, 103.617063048169, 159073532),
    @(31, 87.40588183275602, 99.44045257568359, 101.3525818445749, 81.77792171181505, 159073532),
    @(32, 75.16998855284349, 80.281982421875, 84.08143559484539, 66.91974350507481, 159073532),
    @(33, 87.78619979798464, 100.4823455810547, 103.7131981790559, 84.52571513332373, 159073532),
    @(34, 105.3317272614105, 90.37758636474609, 106.1328393921542, 87.56874392637896, 159073532),
    @(35, 110.6584657251123, 111.8268356323242, 118.0251299624344, 105.0344490190034, 159073532),
    @(36, 99.62141071326366, 82.53382110595703, 100.7315041486907, 81.22549129852919, 159073532),
    @(37, 106.180367060856, 95.85010528564452, 112.7298087485501, 91.58304834820029, 159073532),
    @(38, 112.1180198492091, 115.5655136108398, 118.0592351273358, 95.17861576377771, 159073532),
    @(39, 147.1666188520728, 130.4214477539062, 162.2909787223239, 119.4237209430777, 159073532),
    @(40, 132.1773449063756, 105.7120132446289, 134.7253394281065, 105.0352012505913, 159073532),
    @(41, 126.3057255923791, 115.3750381469727, 143.6433698097777, 100.4088618605841, 159073532),
    @(42, 81.67341275865937, 74.02276611328125, 84.67582061311903, 65.60405792808659, 159073532),
    @(43, 89.43001932157195, 107.322021484375, 113.3159908325742, 88.51094585072903, 159073532)
)

foreach ($entry in $rowData) {
    $row = $entry[0]
    $ws.Cells.Item($row, 4).Value = $entry[1]   # open_price
    $ws.Cells.Item($row, 5).Value = $entry[2]   # close_price
    $ws.Cells.Item($row, 6).Value = $entry[3]   # high_price
    $ws.Cells.Item($row, 7).Value = $entry[4]   # low_price
    $ws.Cells.Item($row, 8).Value = $entry[5]   # shares_outstanding
    $ws.Cells.Item($row, 9).Value = "TER"       # fixed_ticker -> TER
}
